$wb = $excel.ActiveWorkbook

# Insert the new "Week 5" sheet after the last existing sheet (becomes the active sheet)
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "Week 5"

# Column widths (closest representable value to the source file's 35.88671875)
$ws.Columns.Item(3).ColumnWidth = 35

# Date cells (numFmtId 16 => "d-mmm", matches style used on the other week sheets)
$ws.Range("B2").NumberFormat = "d-mmm"
$ws.Range("B2").Value = 42793

$ws.Range("B4").NumberFormat = "d-mmm"
$ws.Range("B4").Value = 42794

$ws.Range("B6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = 42795

$ws.Range("B7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = 42796

# Text/progress entries (order matches the shared-string insertion order of the source file)
$ws.Range("C2").Value = "Project poster hand in/presentation"
$ws.Range("C3").Value = "started work on report"
$ws.Range("C4").Value = "created basic ui"
$ws.Range("C5").Value = "compass implemented, reading values from the sensors"
$ws.Range("C6").Value = "splitting received locations into north, south, east and west"
$ws.Range("G3").Value = "possible battery drain issue"
$ws.Range("C7").Value = "work on report"

# Selection/active cell on the new sheet
$ws.Range("C8").Select() | Out-Null
